$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff/Handback Datetime for the second data row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-20 20:39:41"
$wsZhCn.Range("H3").Value = "2016-03-20 20:40:01"

# de-de sheet: update Correspond Handoff/Handback Datetime for the second data row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-20 20:39:44"
$wsDeDe.Range("H3").Value = "2016-03-20 20:40:08"
